# SAV-700: Update charts importer test fixture.
# - Ensure the ID of the "charting date" recorded program data element is
#   constant ("testchartcode0" -> "PatientChartingDate" on the Test Chart
#   sheet).
# - Ensure the ID of all complex chart CORE questions is constant
#   ("testchartcorecode0".."testchartcorecode4" -> the already-used
#   human-readable question codes on the Core sheet).
# - Widen column A on the Test Chart sheet to fit the new, longer code.

$wb = $excel.ActiveWorkbook
$core = $wb.Worksheets.Item("Core")
$chart = $wb.Worksheets.Item("Test Chart")

# ---------------------------------------------------------------------
# Grab the distinguished formatting (the "s=2" style) that currently
# lives on Core!A6 BEFORE that row gets normalised below, and stamp it
# onto Test Chart!A2 (which adopts that same distinguished style in the
# updated fixture).
# ---------------------------------------------------------------------
$core.Range("A6").Copy() | Out-Null
$chart.Range("A2").PasteSpecial(-4122) | Out-Null
$chart.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# "Core" sheet: collapse the per-row placeholder question codes in
# column A onto the real (already reused in column B) question codes.
# ---------------------------------------------------------------------
$core.Range("A2").Value = "ComplexChartInstanceName"
$core.Range("A3").Value = "ComplexChartDate"
$core.Range("A4").Value = "ComplexChartType"
$core.Range("A5").Value = "ComplexChartSubtype"
$core.Range("A6").Value = "ComplexChartInstanceName"

# Row 6 (A6/C6/D6) previously used the alternate "s=2" style; the updated
# fixture normalises it to the same style as every other data row ("s=1"),
# matching row 2's A/B cells. Copy formatting only (keep the values/shared
# strings untouched) from a cell that already has the target style.
$core.Range("B6").Copy() | Out-Null
$core.Range("A6:D6").PasteSpecial(-4122) | Out-Null
$core.Application.CutCopyMode = $false
$core.Range("A6").Value = "ComplexChartInstanceName"
$core.Range("B6").Value = "ComplexChartInstanceName"

# ---------------------------------------------------------------------
# "Test Chart" sheet: rename the first data element's code and widen
# column A so the longer code is fully visible.
# ---------------------------------------------------------------------
$chart.Range("A2").Value = "PatientChartingDate"

# NOTE: Excel's ColumnWidth (character units) is re-derived from pixel width
# on save, which is not a 1:1 mapping with the stored XML `width` (~0.83
# wider). 16.67 characters round-trips to exactly width="17.5" in the saved
# OOXML, matching the target column width.
$chart.Columns.Item(1).ColumnWidth = 16.67
